# Auto-generated Excel COM-interop script
# Updates computed market-price / profit columns (H:N) across all 8 class sheets
# to refreshed values pulled by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# ALC: set updated values
$ALC_updates = @{
    "H42" = 238.55556
    "I42" = 29.714285
    "J42" = 969.5
    "K42" = 89.142855
    "L42" = 2908.5
    "M42" = 140.857145
    "N42" = -3368.5
    "H62" = 4916.1665
    "J62" = 4916.1665
    "L62" = 4916.1665
    "N62" = -6164.1665
    "H65" = 4916.1665
    "J65" = 4916.1665
    "L65" = 24580.8325
    "N65" = -30820.8325
    "H86" = 250002000
    "I86" = 333334850
    "J86" = 3500
    "K86" = 333334850
    "L86" = 3500
    "M86" = -333333727
    "N86" = -5746
    "H89" = 250002000
    "I89" = 333334850
    "J89" = 3500
    "K89" = 1666674250
    "L89" = 17500
    "M89" = -1666668634
    "N89" = -28732
    "H101" = 3002.2307
    "I101" = 406.14285
    "J101" = 6031
    "K101" = 1218.42855
    "L101" = 18093
    "M101" = 403.5714499999999
    "N101" = -21337
    "H106" = 3649.3333
    "I106" = 3000
    "K106" = 3000
    "M106" = -2369
    "H132" = 4266.3945
    "I132" = 1170.8788
    "J132" = 24696.8
    "K132" = 3512.6364
    "L132" = 74090.39999999999
    "M132" = -982.6363999999999
    "N132" = -79150.39999999999
    "H137" = 9526372
    "I137" = 2052.2307
    "K137" = 6156.6921
    "M137" = -3606.6921
}
foreach ($cell in $ALC_updates.Keys) {
    $ws.Range($cell).Value = $ALC_updates[$cell]
}

$ws = $wb.Worksheets.Item("ARM")

# ARM: set updated values
$ARM_updates = @{
    "H32" = 14465.193
    "I32" = 13787.389
    "K32" = 13787.389
    "M32" = -13500.389
    "H61" = 4800.7856
    "I61" = 2837.3635
    "K61" = 2837.3635
    "M61" = -2625.3635
    "H104" = 45056
    "J104" = 45056
    "L104" = 45056
    "N104" = -52044
    "H110" = 2318.5557
    "I110" = 2271.2856
    "J110" = 2484
    "K110" = 2271.2856
    "L110" = 2484
    "M110" = -226.2856000000002
    "N110" = -6574
    "H132" = 1465.8334
    "I132" = 1465.8334
    "K132" = 4397.5002
    "M132" = -1867.5002
    "H136" = 4800.7856
    "I136" = 2837.3635
    "K136" = 8512.0905
    "M136" = -5962.0905
    "H139" = 81579.39999999999
    "J139" = 81579.39999999999
    "L139" = 81579.39999999999
    "N139" = -91859.39999999999
}
foreach ($cell in $ARM_updates.Keys) {
    $ws.Range($cell).Value = $ARM_updates[$cell]
}

$ws = $wb.Worksheets.Item("BSM")

# BSM: set updated values
$BSM_updates = @{
    "H20" = 4723.722
    "J20" = 5977.125
    "L20" = 5977.125
    "N20" = -6471.125
    "H22" = 899.63635
    "I22" = 789.4
    "K22" = 789.4
    "M22" = -616.4
    "H105" = 5378.6924
    "I105" = 5124.1113
    "J105" = 5951.5
    "K105" = 5124.1113
    "L105" = 5951.5
    "M105" = -3377.1113
    "N105" = -9445.5
}
foreach ($cell in $BSM_updates.Keys) {
    $ws.Range($cell).Value = $BSM_updates[$cell]
}

$ws = $wb.Worksheets.Item("CRP")

# CRP: set updated values
$CRP_updates = @{
    "J31" = 5761.4443
    "L31" = 5761.4443
    "N31" = -6351.4443
    "J34" = 5761.4443
    "L34" = 5761.4443
    "N34" = -6165.4443
    "H58" = 2179.5881
    "I58" = 2058.6875
    "J58" = 4114
    "K58" = 2058.6875
    "L58" = 4114
    "M58" = -1855.6875
    "N58" = -4520
    "H94" = 2219.7273
    "I94" = 1197.5
    "J94" = 2446.889
    "K94" = 1197.5
    "L94" = 2446.889
    "M94" = -746.5
    "N94" = -3348.889
    "H107" = 943.9524
    "J107" = 1060.6666
    "L107" = 1060.6666
    "N107" = -4900.6666
    "H132" = 88896940
    "I132" = 95239580
    "J132" = 100000
    "K132" = 285718740
    "L132" = 300000
    "M132" = -285716210
    "N132" = -305060
    "H136" = 2179.5881
    "I136" = 2058.6875
    "J136" = 4114
    "K136" = 6176.0625
    "L136" = 12342
    "M136" = -3626.0625
    "N136" = -17442
}
foreach ($cell in $CRP_updates.Keys) {
    $ws.Range($cell).Value = $CRP_updates[$cell]
}

$ws = $wb.Worksheets.Item("CUL")

# CUL: set updated values
$CUL_updates = @{
    "H4" = 7192027.5
    "I4" = 316713
    "K4" = 950139
    "M4" = -950027
    "H5" = 615.2353000000001
    "I5" = 180.25
    "K5" = 540.75
    "M5" = -428.75
    "H12" = 173.61111
    "J12" = 139
    "L12" = 417
    "N12" = -763
    "H100" = 3010
    "J100" = 3010
    "L100" = 9030
    "N100" = -10652
    "H122" = 1661.4286
    "I122" = 0
    "K122" = 0
    "H135" = 615.2353000000001
    "I135" = 180.25
    "K135" = 1622.25
    "M135" = 912.75
}
foreach ($cell in $CUL_updates.Keys) {
    $ws.Range($cell).Value = $CUL_updates[$cell]
}

# CUL: clear cells removed in the update
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")

# GSM: set updated values
$GSM_updates = @{
    "H35" = 0
    "J35" = 0
    "L35" = 0
    "H80" = 80257.47
    "I80" = 115279.664
    "K80" = 115279.664
    "M80" = -114281.664
    "H83" = 80257.47
    "I83" = 115279.664
    "K83" = 576398.3200000001
    "M83" = -571406.3200000001
    "H111" = 40000
    "J111" = 40000
    "L111" = 40000
    "N111" = -46134
    "H113" = 4383
    "I113" = 1574.5
    "J113" = 10000
    "K113" = 1574.5
    "L113" = 10000
    "M113" = 595.5
    "N113" = -14340
    "H126" = 2990.6538
    "I126" = 2178.8333
    "J126" = 3686.5
    "K126" = 6536.499899999999
    "L126" = 11059.5
    "M126" = -4066.499899999999
    "N126" = -15999.5
    "H132" = 2688
    "J132" = 6225
    "L132" = 18675
    "N132" = -23735
}
foreach ($cell in $GSM_updates.Keys) {
    $ws.Range($cell).Value = $GSM_updates[$cell]
}

# GSM: clear cells removed in the update
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("LTW")

# LTW: set updated values
$LTW_updates = @{
    "H55" = 302.78946
    "I55" = 112
    "J55" = 565.125
    "K55" = 112
    "L55" = 565.125
    "M55" = 61
    "N55" = -911.125
    "H122" = 9058.392
    "I122" = 4580
    "J122" = 13163.583
    "K122" = 13740
    "L122" = 39490.749
    "M122" = -11290
    "N122" = -44390.749
    "H136" = 4096.952
    "I136" = 2806.1714
    "J136" = 5710.4287
    "K136" = 8418.514200000001
    "L136" = 17131.2861
    "M136" = -5868.514200000001
    "N136" = -22231.2861
    "H139" = 86712
    "J139" = 86712
    "L139" = 86712
    "N139" = -96992
}
foreach ($cell in $LTW_updates.Keys) {
    $ws.Range($cell).Value = $LTW_updates[$cell]
}

$ws = $wb.Worksheets.Item("WVR")

# WVR: set updated values
$WVR_updates = @{
    "H32" = 19166.666
    "I32" = 7500
    "K32" = 7500
    "M32" = -7183
    "H107" = 783.41174
    "I107" = 685.375
    "J107" = 870.55554
    "K107" = 2056.125
    "L107" = 2611.66662
    "M107" = -136.125
    "N107" = -6451.66662
    "H113" = 364.89655
    "I113" = 354.64706
    "J113" = 379.41666
    "K113" = 1063.94118
    "L113" = 1138.24998
    "M113" = 1106.05882
    "N113" = -5478.249980000001
    "H132" = 22224840
    "I132" = 30304196
    "K132" = 90912588
    "M132" = -90910058
    "H133" = 104999.5
    "J133" = 104999.5
    "L133" = 104999.5
    "N133" = -115119.5
    "H136" = 2340.123
    "I136" = 1652.46
    "K136" = 4957.38
    "M136" = -2407.38
}
foreach ($cell in $WVR_updates.Keys) {
    $ws.Range($cell).Value = $WVR_updates[$cell]
}
